$d = $word.ActiveDocument

function Find-ParaIndex($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# --- "Milestone 2" body ("Cross session ...") gets marked complete like Milestone 1,
#     i.e. prefixed with a new "X - " run. ---
$idxCross = Find-ParaIndex $d "Cross session*"
$pCross = $d.Paragraphs.Item($idxCross)
$pCross.Range.InsertBefore("X - ")

# --- New bullet inserted right after "Milestone 3:" and before "Finalize UI panel...":
#     "Warning for automatically resuming with don't remind me again option" ---
$idxM3 = Find-ParaIndex $d "*Milestone 3:*"
$pM3 = $d.Paragraphs.Item($idxM3)
$pM3.Range.InsertParagraphAfter()
$pWarning = $d.Paragraphs.Item($idxM3 + 1)
$pWarning.Range.Text = "`t`tX - Warning for automatically resuming with don" + [char]8217 + "t remind me again option"

# --- "Indicator icon in status bar" gets marked complete too: insert "X - " right
#     before the text, after its two leading tabs. ---
$idxIndicator = Find-ParaIndex $d "*Indicator icon in status bar*"
$pIndicator = $d.Paragraphs.Item($idxIndicator)
$fullIndicator = $pIndicator.Range
$insertAt = $fullIndicator.Start + 2
$insertRange = $d.Range($insertAt, $insertAt)
$insertRange.InsertBefore("X - ")

# --- Two new bullets appended after "Allow name format specification":
#     "Allow file format specification and compression" and "Preferences in addon settings" ---
$idxAllowName = Find-ParaIndex $d "*Allow name format specification*"
$pAllowName = $d.Paragraphs.Item($idxAllowName)
$pAllowName.Range.InsertParagraphAfter()
$pFileFormat = $d.Paragraphs.Item($idxAllowName + 1)
$pFileFormat.Range.Text = "`t`tAllow file format specification and compression"

$pFileFormat.Range.InsertParagraphAfter()
$pPreferences = $d.Paragraphs.Item($idxAllowName + 2)
$pPreferences.Range.Text = "`t`tPreferences in addon settings"

Write-Output "done"
